# Add a new "time_taken" column (F) with a header styled like the other
# header cells, and a timestamp value for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same formatting as the other header cells (copy from E1)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Timestamp values for rows 2-15
$timestamps = @(
    "2021-10-05 13:41:14.585714",
    "2021-10-05 13:41:14.585725",
    "2021-10-05 13:41:14.585729",
    "2021-10-05 13:41:14.585733",
    "2021-10-05 13:41:14.585736",
    "2021-10-05 13:41:14.585739",
    "2021-10-05 13:41:14.585743",
    "2021-10-05 13:41:14.585746",
    "2021-10-05 13:41:14.585749",
    "2021-10-05 13:41:14.585752",
    "2021-10-05 13:41:14.585755",
    "2021-10-05 13:41:14.585758",
    "2021-10-05 13:41:14.585761",
    "2021-10-05 13:41:14.585764"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

$excel.CutCopyMode = 0
